# Natmi following Dr Hou advice
# Update existing rows 2-4 with recalculated values (ligand-expressing cells count
# changed from 1 to 3 across the dataset, new totals/averages etc.), and append a
# new row 5 for the additional sending-cluster pairing (sCs -> M2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> M2) ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.637903
$ws.Range("H2").Value = 1.913709
$ws.Range("I2").Value = 0.1229013127714845
$ws.Range("J2").Value = 0.1229013127714844
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.557669333333333
$ws.Range("N2").Value = 28.673008
$ws.Range("Q2").Value = 6.096865940741333
$ws.Range("R2").Value = 54.871793466672
$ws.Range("S2").Value = 0.1229013127714845
$ws.Range("T2").Value = 0.1229013127714844

# --- Row 3 (FAPs -> M2) ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.170281
$ws.Range("H3").Value = 0.510843
$ws.Range("I3").Value = 0.03280711713229307
$ws.Range("J3").Value = 0.03280711713229307
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.557669333333333
$ws.Range("N3").Value = 28.673008
$ws.Range("Q3").Value = 1.627489491749333
$ws.Range("R3").Value = 14.647405425744
$ws.Range("S3").Value = 0.03280711713229307
$ws.Range("T3").Value = 0.03280711713229307

# --- Row 4 (M2 -> M2) ---
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.407124666666667
$ws.Range("H4").Value = 10.221374
$ws.Range("I4").Value = 0.6564322386153376
$ws.Range("J4").Value = 0.6564322386153377
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.557669333333333
$ws.Range("N4").Value = 28.673008
$ws.Range("Q4").Value = 32.56417094144356
$ws.Range("R4").Value = 293.077538472992
$ws.Range("S4").Value = 0.6564322386153376
$ws.Range("T4").Value = 0.6564322386153377

# --- Row 5 (new sending cluster: sCs -> M2) ---
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Tnfsf13"
$ws.Range("C5").Value = "Tnfrsf13b"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.975059
$ws.Range("H5").Value = 2.925177
$ws.Range("I5").Value = 0.1878593314808848
$ws.Range("J5").Value = 0.1878593314808848
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.557669333333333
$ws.Range("N5").Value = 28.673008
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 9.319291502490666
$ws.Range("R5").Value = 83.873623522416
$ws.Range("S5").Value = 0.1878593314808848
$ws.Range("T5").Value = 0.1878593314808848
